$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated GDP values for existing rows (B12:B30)
$updates = @{
    12 = 242677.43700000001
    13 = 249723.97899999999
    14 = 264853.37099999998
    15 = 282314.65399999998
    16 = 293826.00199999998
    17 = 306907.95299999998
    18 = 326976.46799999999
    19 = 340107.723
    20 = 334234.96399999998
    21 = 343539.46500000003
    22 = 351690.65899999999
    23 = 367778.658
    24 = 377541.28899999999
    25 = 394728.109
    26 = 407824.14500000002
    27 = 414883.70600000001
    28 = 418695.98100000003
    29 = 435423.85200000001
    30 = 450742.19199999998
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}

# Add new row 31 for 2020-01-01 observation
$ws.Cells.Item(31, 1).Value = 43831
$ws.Cells.Item(31, 1).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(31, 2).Value = 439055.11800000002
$ws.Cells.Item(31, 2).NumberFormat = "0.000"

# Update the selection to match the new state (select full columns A:B)
$ws.Range("A:B").Select() | Out-Null
